$d = $word.ActiveDocument

$d.Content.Find.Execute("53×23=", $true, $false, $false, $false, $false, $true, 1, $false, "32×68=", 2) | Out-Null
$d.Content.Find.Execute("95×70=", $true, $false, $false, $false, $false, $true, 1, $false, "75×37=", 2) | Out-Null
$d.Content.Find.Execute("58×64=", $true, $false, $false, $false, $false, $true, 1, $false, "58×85=", 2) | Out-Null
$d.Content.Find.Execute("86×44=", $true, $false, $false, $false, $false, $true, 1, $false, "56×16=", 2) | Out-Null
$d.Content.Find.Execute("85×17=", $true, $false, $false, $false, $false, $true, 1, $false, "88×85=", 2) | Out-Null
$d.Content.Find.Execute("25×87=", $true, $false, $false, $false, $false, $true, 1, $false, "34×63=", 2) | Out-Null
$d.Content.Find.Execute("33×55=", $true, $false, $false, $false, $false, $true, 1, $false, "68×54=", 2) | Out-Null
$d.Content.Find.Execute("99×18=", $true, $false, $false, $false, $false, $true, 1, $false, "14×71=", 2) | Out-Null
$d.Content.Find.Execute("93×99=", $true, $false, $false, $false, $false, $true, 1, $false, "81×54=", 2) | Out-Null
$d.Content.Find.Execute("15×38=", $true, $false, $false, $false, $false, $true, 1, $false, "34×70=", 2) | Out-Null
$d.Content.Find.Execute("75×30=", $true, $false, $false, $false, $false, $true, 1, $false, "18×23=", 2) | Out-Null
$d.Content.Find.Execute("46×83=", $true, $false, $false, $false, $false, $true, 1, $false, "61×56=", 2) | Out-Null
$d.Content.Find.Execute("75×63=", $true, $false, $false, $false, $false, $true, 1, $false, "80×58=", 2) | Out-Null
$d.Content.Find.Execute("59×60=", $true, $false, $false, $false, $false, $true, 1, $false, "68×37=", 2) | Out-Null
$d.Content.Find.Execute("70×82=", $true, $false, $false, $false, $false, $true, 1, $false, "69×69=", 2) | Out-Null
$d.Content.Find.Execute("18×57=", $true, $false, $false, $false, $false, $true, 1, $false, "22×95=", 2) | Out-Null
$d.Content.Find.Execute("61×76=", $true, $false, $false, $false, $false, $true, 1, $false, "79×57=", 2) | Out-Null
$d.Content.Find.Execute("28×17=", $true, $false, $false, $false, $false, $true, 1, $false, "82×43=", 2) | Out-Null
$d.Content.Find.Execute("40×51=", $true, $false, $false, $false, $false, $true, 1, $false, "38×78=", 2) | Out-Null
$d.Content.Find.Execute("61×66=", $true, $false, $false, $false, $false, $true, 1, $false, "62×91=", 2) | Out-Null
$d.Content.Find.Execute("58×40=", $true, $false, $false, $false, $false, $true, 1, $false, "73×57=", 2) | Out-Null
$d.Content.Find.Execute("79×70=", $true, $false, $false, $false, $false, $true, 1, $false, "53×56=", 2) | Out-Null
$d.Content.Find.Execute("36×79=", $true, $false, $false, $false, $false, $true, 1, $false, "73×59=", 2) | Out-Null
$d.Content.Find.Execute("52×49=", $true, $false, $false, $false, $false, $true, 1, $false, "17×19=", 2) | Out-Null
$d.Content.Find.Execute("13×54=", $true, $false, $false, $false, $false, $true, 1, $false, "57×23=", 2) | Out-Null
